$wb = $excel.ActiveWorkbook

# zh-cn sheet: Correspond Handoff Datetime (E) and Correspond Handback DateTime (H)
# Both rows 2 and 3 shared the same datetime text in the original file, so update both.
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-24 11:00:26"
$wsZh.Range("E3").Value = "2016-03-24 11:00:26"
$wsZh.Range("H2").Value = "2016-03-24 11:01:43"
$wsZh.Range("H3").Value = "2016-03-24 11:01:43"

# de-de sheet: Correspond Handoff Datetime (E) and Correspond Handback DateTime (H)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-24 11:00:40"
$wsDe.Range("E3").Value = "2016-03-24 11:00:40"
$wsDe.Range("H2").Value = "2016-03-24 11:02:02"
$wsDe.Range("H3").Value = "2016-03-24 11:02:02"
